# fix: some bugs and ignore some item not calculate
#
# The 薪資表 (Salary table) worksheet (the first / active sheet) had a value
# in column S of the first data row (row 5) that should not have been
# included in the deduction calculation ("ignore some item not calculate").
# Restore the correct amount of 300. The dependent formulas on the same row
# (X5 = SUM(Q5:W5), Y5 = P5-X5, AF5 = Y5+AE5) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("S5").Value = 300

# Leave the selection on the cell the user ended up on after the edit.
$ws.Range("S6").Select()
